$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K6").Value = 754.71429
$ws.Range("H6").Value = 1383.6471
$ws.Range("I6").Value = 251.57143
$ws.Range("M6").Value = -642.71429
$ws.Range("I15").Value = 2265.1392
$ws.Range("H15").Value = 2265.1392
$ws.Range("K15").Value = 6795.417600000001
$ws.Range("M15").Value = -6626.417600000001
$ws.Range("J17").Value = 20345.785
$ws.Range("L17").Value = 61037.355
$ws.Range("N17").Value = -61373.355
$ws.Range("H17").Value = 20345.785
$ws.Range("K31").Value = 5073999.300000001
$ws.Range("M31").Value = -5073769.300000001
$ws.Range("I31").Value = 1691333.1
$ws.Range("H31").Value = 1691333.1
$ws.Range("M94").Value = -4910
$ws.Range("K94").Value = 5361
$ws.Range("H94").Value = 6288.8
$ws.Range("I94").Value = 5361
$ws.Range("M99").Value = 181
$ws.Range("H99").Value = 451.2857
$ws.Range("I99").Value = 439
$ws.Range("L99").Value = 1489.00002
$ws.Range("K99").Value = 1317
$ws.Range("J99").Value = 496.33334
$ws.Range("N99").Value = -4485.000019999999
$ws.Range("I100").Value = 1682.091
$ws.Range("L100").Value = 2971.0715
$ws.Range("N100").Value = -4053.0715
$ws.Range("K100").Value = 1682.091
$ws.Range("J100").Value = 2971.0715
$ws.Range("M100").Value = -1141.091
$ws.Range("H100").Value = 2403.92
$ws.Range("M135").Value = -1597.28565
$ws.Range("J135").Value = 1382.5
$ws.Range("I135").Value = 459.14285
$ws.Range("L135").Value = 12442.5
$ws.Range("N135").Value = -17512.5
$ws.Range("K135").Value = 4132.28565
$ws.Range("H135").Value = 574.5625
$ws.Range("M138").Value = 1082.5
$ws.Range("K138").Value = 4057.5
$ws.Range("J138").Value = 3991.4736
$ws.Range("N138").Value = -22254.4208
$ws.Range("H138").Value = 2355.31
$ws.Range("L138").Value = 11974.4208
$ws.Range("I138").Value = 1352.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J26").Value = 20000
$ws.Range("I26").Value = 2989.4
$ws.Range("L26").Value = 20000
$ws.Range("K26").Value = 2989.4
$ws.Range("H26").Value = 5824.5
$ws.Range("M26").Value = -2659.4
$ws.Range("N26").Value = -20660
$ws.Range("K45").Value = 1080.75
$ws.Range("J45").Value = 1050
$ws.Range("L45").Value = 1050
$ws.Range("I45").Value = 1080.75
$ws.Range("H45").Value = 1070.5
$ws.Range("M45").Value = -703.75
$ws.Range("N45").Value = -1804
$ws.Range("L97").Value = 913.75
$ws.Range("K97").Value = 603.5833
$ws.Range("M97").Value = -107.5833
$ws.Range("H97").Value = 681.125
$ws.Range("N97").Value = -1905.75
$ws.Range("J97").Value = 913.75
$ws.Range("I97").Value = 603.5833
$ws.Range("J101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("J102").Value = 334433
$ws.Range("N102").Value = -337677
$ws.Range("L102").Value = 334433
$ws.Range("H102").Value = 46696.727
$ws.Range("K102").Value = 1264.6842
$ws.Range("I102").Value = 1264.6842
$ws.Range("M102").Value = 357.3158000000001
$ws.Range("N122").Value = -10161.25
$ws.Range("I122").Value = 2300
$ws.Range("H122").Value = 1935.8334
$ws.Range("L122").Value = 5261.25
$ws.Range("M122").Value = -4450
$ws.Range("K122").Value = 6900
$ws.Range("J122").Value = 1753.75

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J33").Value = 6000
$ws.Range("I33").Value = 0
$ws.Range("H33").Value = 6000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 6000
$ws.Range("N33").Value = -6672
$ws.Range("M33").ClearContents()
$ws.Range("L38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("L79").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("J94").Value = 203000
$ws.Range("L94").Value = 203000
$ws.Range("M94").Value = -1304.1111
$ws.Range("K94").Value = 1755.1111
$ws.Range("H94").Value = 73628.28999999999
$ws.Range("I94").Value = 1755.1111
$ws.Range("N94").Value = -203902
$ws.Range("M99").Value = 639.25
$ws.Range("H99").Value = 1484.8
$ws.Range("I99").Value = 858.75
$ws.Range("L99").Value = 3989
$ws.Range("K99").Value = 858.75
$ws.Range("J99").Value = 3989
$ws.Range("N99").Value = -6985

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L20").Value = 25000
$ws.Range("H20").Value = 25000
$ws.Range("J20").Value = 25000
$ws.Range("N20").Value = -25472
$ws.Range("J30").Value = 25000
$ws.Range("L30").Value = 25000
$ws.Range("H30").Value = 25000
$ws.Range("N30").Value = -25182
$ws.Range("I32").Value = 2950
$ws.Range("H32").Value = 2950
$ws.Range("K32").Value = 2950
$ws.Range("M32").Value = -2634
$ws.Range("H38").Value = 3200
$ws.Range("I38").Value = 3200
$ws.Range("K38").Value = 3200
$ws.Range("M38").Value = -2823
$ws.Range("H46").Value = 3200
$ws.Range("M46").Value = -2989
$ws.Range("I46").Value = 3200
$ws.Range("K46").Value = 3200
$ws.Range("H58").Value = 976084.8
$ws.Range("J58").Value = 2134
$ws.Range("N58").Value = -2540
$ws.Range("I58").Value = 1544222.9
$ws.Range("L58").Value = 2134
$ws.Range("K58").Value = 1544222.9
$ws.Range("M58").Value = -1544019.9
$ws.Range("M99").Value = -3218
$ws.Range("H99").Value = 3900.8572
$ws.Range("I99").Value = 4716
$ws.Range("L99").Value = 3085.7144
$ws.Range("K99").Value = 4716
$ws.Range("J99").Value = 3085.7144
$ws.Range("N99").Value = -6081.7144
$ws.Range("N122").Value = -11790.4999
$ws.Range("I122").Value = 3371.0667
$ws.Range("H122").Value = 3064.1428
$ws.Range("L122").Value = 6890.499899999999
$ws.Range("M122").Value = -7663.2001
$ws.Range("K122").Value = 10113.2001
$ws.Range("J122").Value = 2296.8333
$ws.Range("L123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("N123").Value = -89800
$ws.Range("H123").Value = 80000
$ws.Range("J126").Value = 3085.7144
$ws.Range("M126").Value = -11678
$ws.Range("H126").Value = 3900.8572
$ws.Range("I126").Value = 4716
$ws.Range("K126").Value = 14148
$ws.Range("L126").Value = 9257.143199999999
$ws.Range("N126").Value = -14197.1432
$ws.Range("H128").Value = 25000
$ws.Range("L128").Value = 25000
$ws.Range("J128").Value = 25000
$ws.Range("N128").Value = -34960
$ws.Range("I136").Value = 1544222.9
$ws.Range("L136").Value = 6402
$ws.Range("K136").Value = 4632668.699999999
$ws.Range("M136").Value = -4630118.699999999
$ws.Range("H136").Value = 976084.8
$ws.Range("N136").Value = -11502
$ws.Range("J136").Value = 2134

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L64").Value = 14593.4208
$ws.Range("K64").Value = 9711.999899999999
$ws.Range("J64").Value = 4864.4736
$ws.Range("M64").Value = -9441.999899999999
$ws.Range("H64").Value = 4642.591
$ws.Range("N64").Value = -15133.4208
$ws.Range("I64").Value = 3237.3333
$ws.Range("M67").Value = -8775.999899999999
$ws.Range("L67").Value = 14593.4208
$ws.Range("H67").Value = 4642.591
$ws.Range("K67").Value = 9711.999899999999
$ws.Range("N67").Value = -16465.4208
$ws.Range("J67").Value = 4864.4736
$ws.Range("I67").Value = 3237.3333
$ws.Range("H131").Value = 2519.42
$ws.Range("N131").Value = -14811.5454
$ws.Range("J131").Value = 1577.1818
$ws.Range("L131").Value = 4731.5454

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K70").Value = 6397.143
$ws.Range("H70").Value = 7649.8335
$ws.Range("N70").Value = -9943.6
$ws.Range("M70").Value = -6127.143
$ws.Range("I70").Value = 6397.143
$ws.Range("J70").Value = 9403.6
$ws.Range("L70").Value = 9403.6
$ws.Range("K73").Value = 6397.143
$ws.Range("N73").Value = -11275.6
$ws.Range("M73").Value = -5461.143
$ws.Range("H73").Value = 7649.8335
$ws.Range("I73").Value = 6397.143
$ws.Range("L73").Value = 9403.6
$ws.Range("J73").Value = 9403.6
$ws.Range("K113").Value = 990
$ws.Range("M113").Value = 1180
$ws.Range("I113").Value = 990
$ws.Range("J113").Value = 2159.6365
$ws.Range("N113").Value = -6499.636500000001
$ws.Range("H113").Value = 1847.7333
$ws.Range("L113").Value = 2159.6365
$ws.Range("N122").Value = -15658
$ws.Range("I122").Value = 1710.3636
$ws.Range("H122").Value = 2112.2856
$ws.Range("L122").Value = 10758
$ws.Range("M122").Value = -2681.0908
$ws.Range("K122").Value = 5131.0908
$ws.Range("J122").Value = 3586
$ws.Range("N132").Value = -13904.3999
$ws.Range("H132").Value = 2171.8708
$ws.Range("K132").Value = 4332.375
$ws.Range("M132").Value = -1802.375
$ws.Range("L132").Value = 8844.3999
$ws.Range("I132").Value = 1444.125
$ws.Range("J132").Value = 2948.1333

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 5503
$ws.Range("L7").Value = 5756.1113
$ws.Range("H7").Value = 5678.231
$ws.Range("N7").Value = -5980.1113
$ws.Range("M7").Value = -5391
$ws.Range("J7").Value = 5756.1113
$ws.Range("K7").Value = 5503
$ws.Range("N40").Value = -3870
$ws.Range("H40").Value = 2475
$ws.Range("L40").Value = 3598
$ws.Range("K40").Value = 1964.5454
$ws.Range("I40").Value = 1964.5454
$ws.Range("M40").Value = -1828.5454
$ws.Range("J40").Value = 3598
$ws.Range("M93").Value = 574.1
$ws.Range("K93").Value = 673.9
$ws.Range("H93").Value = 673.9
$ws.Range("I93").Value = 673.9
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("J126").Value = 5756.1113
$ws.Range("M126").Value = -14039
$ws.Range("H126").Value = 5678.231
$ws.Range("I126").Value = 5503
$ws.Range("K126").Value = 16509
$ws.Range("N126").Value = -22208.3339
$ws.Range("L126").Value = 17268.3339

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 500000
$ws.Range("L99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("N132").Value = -9032.428400000001
$ws.Range("H132").Value = 1108.0454
$ws.Range("K132").Value = 3021.6
$ws.Range("M132").Value = -491.6000000000004
$ws.Range("L132").Value = 3972.4284
$ws.Range("I132").Value = 1007.2
$ws.Range("J132").Value = 1324.1428
$ws.Range("I136").Value = 976.9545000000001
$ws.Range("L136").Value = 3669.75
$ws.Range("K136").Value = 2930.8635
$ws.Range("M136").Value = -380.8635000000004
$ws.Range("H136").Value = 1014.8461
$ws.Range("N136").Value = -8769.75
$ws.Range("J136").Value = 1223.25
